$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 249. This shifts the existing data
# block (rows 249:354) down to (250:355), preserving all of its values
# and formatting, and leaves row 249 empty and ready for the new
# weekly data point.
$ws.Rows("249:249").Insert()

# Populate the newly inserted row 249 with the new weekly record.
$ws.Range("A249").Value = 10
$ws.Range("B249").Value = "Vega Modelo de Temuco"
$ws.Range("C249").Value = "La Araucanía"
$ws.Range("D249").Value = 44523
$ws.Range("E249").Value = 9
$ws.Range("F249").Value = 100112043
$ws.Range("G249").Value = "Pepino ensalada"
$ws.Range("H249").Value = "Sin especificar"
$ws.Range("I249").Value = "Primera"
$ws.Range("J249").Value = 250
$ws.Range("K249").Value = 8000
$ws.Range("L249").Value = 10000
$ws.Range("M249").Value = 8800
$ws.Range("N249").Value = "$/caja 60 unidades"
$ws.Range("O249").Value = "Región del Maule"
$ws.Range("P249").Value = 147
$ws.Range("Q249").Value = 60
$ws.Range("R249").Value = "Hortaliza"

# Match the date number format used by the rest of the column.
$ws.Range("D249").NumberFormat = "YYYY-MM-DD HH:MM:SS"
